# Adds test results for "Signavio Process Editor" (three new columns: Import /
# Export / Roundtrip, mirroring the existing "MID Innovator" block) and a new
# vertical-pool test-case row to the "Test Case Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Case Results")

# ---------------------------------------------------------------------------
# 1. Insert a new row for the "vertical pool" test case, right below the
#    existing 4-row "Call Activity - Expanded" block (old row 9, now row 10).
# ---------------------------------------------------------------------------
$ws.Rows("9:9").Insert()

# Re-merge column A for the "Call Activity - Expanded" group (now A5:A9) and
# the "Call Activity - Collapsed" group, which Excel already shifted down to
# A10:A13 automatically.
$ws.Range("A5:A8").UnMerge()
$ws.Range("A5:A9").Merge()

# Give the new row the same look (borders/fills/fonts) as the row above it,
# then move the thick bottom border down from row 8 to row 9 (row 9 is now
# the bottom edge of the merged "Call Activity - Expanded" block).
$ws.Range("A8:K8").Copy()
$ws.Range("A9").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

$ws.Range("A8:K8").Borders.Item(9).LineStyle = -4142
$ws.Range("A9:N9").Borders.Item(9).LineStyle = 1
$ws.Range("A9:N9").Borders.Item(9).Weight = -4138

$ws.Range("C9").Value = 41278
$ws.Range("C9").NumberFormat = "m/d/yyyy"

# ---------------------------------------------------------------------------
# 2. Add the "Signavio Process Editor" result columns (L:N), mirroring the
#    layout used for the existing "MID Innovator" columns (I:K).
# ---------------------------------------------------------------------------
$ws.Range("I2:K13").Copy()
$ws.Range("L2").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

$ws.Range("L2").Value = "Import"
$ws.Range("M2").Value = "Export"
$ws.Range("N2").Value = "Roundtrip"

# Trailing spacer cell matching the header row's thick bottom rule.
$ws.Range("O2").Borders.Item(9).LineStyle = 1
$ws.Range("O2").Borders.Item(9).Weight = -4138

$ws.Range("L3:N3").Merge()
$ws.Range("L3").Value = "Signavio Process Editor"

$ws.Range("L4:N4").Merge()
$ws.Range("L4").Value = 38539

$ws.Range("L5").Value = "OK"
$ws.Range("M5").Value = "OK"
$ws.Range("N5").Value = "OK"

$ws.Range("L6").Value = "OK"
$ws.Range("M6").Value = "OK"
$ws.Range("N6").Value = "OK"

$ws.Range("L7").Value = "OK"
$ws.Range("M7").Value = "OK"
$ws.Range("N7").Value = "OK"

$ws.Range("L8").Value = "Issues"
$ws.Range("M8").Value = "Issues"
$ws.Range("N8").Value = "Issues"

$ws.Range("L9").Value = ""
$ws.Range("M9").Value = ""
$ws.Range("N9").Value = ""

$ws.Range("L10").Value = "Issues"
$ws.Range("M10").Value = "Issues"
$ws.Range("N10").Value = "Issues"

$ws.Range("L11").Value = "Issues"
$ws.Range("M11").Value = "Issues"
$ws.Range("N11").Value = "Issues"

$ws.Range("L12").Value = ""
$ws.Range("M12").Value = ""
$ws.Range("N12").Value = ""

$ws.Range("L13").Value = ""
$ws.Range("M13").Value = ""
$ws.Range("N13").Value = ""

$ws.Range("A2").Select()
